# Add two new trailing columns, I ("I0") and J ("IF"), to the existing
# sheet (previously A1:H76 -> now A1:J76), matching the bold/bordered
# header style already used by column H, and populate the per-row data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells: copy H1's formatting (bold font, thin border,
# centered/top aligned) onto I1:J1, then set their text.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for rows 2-76 (column I = "I0", column J = "IF").
$iVals = @(9,9,9,9,9,10,9,9,9,9,9,9,9,8,9,8,8,7,10,9,9,8,8,9,9,8,9,9,8,9,9,9,9,8,9,9,8,9,9,8,9,9,9,8,8,9,8,8,9,9,8,9,9,9,9,9,9,10,9,8,11,9,8,9,9,8,8,6,5,5,8,7,4,3,3)
$jVals = @(9,9,9,9,9,10,9,9,9,9,9,9,9,9,9,9,9,7,10,9,9,9,9,9,9,8,9,9,9,9,9,9,9,9,9,9,9,10,10,9,9,9,9,9,9,9,9,9,9,9,8,9,9,9,9,9,9,10,9,8,11,9,9,9,9,8,8,6,5,5,8,7,4,3,3)

for ($idx = 0; $idx -lt $iVals.Count; $idx++) {
    $r = 2 + $idx
    $ws.Cells.Item($r, 9).Value = $iVals[$idx]
    $ws.Cells.Item($r, 10).Value = $jVals[$idx]
}

Write-Output "done"
